# "updated chrome driver last version"
# AddCustomerTest!E4 runmode flips from "N" to "Y" (last remaining "N" test
# row is now run), which drops "N" out of the shared-strings table entirely
# once nothing references it anymore.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("AddCustomerTest")
$ws3 = $wb.Worksheets.Item("test_suite")

# Flip the runmode flag for the Bill Gates / AddCustomerTest row.
$ws1.Range("E4").Value = "Y"

# Move the saved selection on test_suite first (selecting there makes it the
# active sheet momentarily) ...
$ws3.Range("B7").Select()

# ... then reselect on AddCustomerTest last so it ends up the active tab,
# matching the workbook's saved state.
$ws1.Range("H5").Select()
